{"js": "// Replace the 25 two-digit-division answers in the single table, in\n// document (row-major) order. Empty \"spacer\" rows are skipped because\n// they contain no text. Each entry is [oldText, newText]; the old text\n// is asserted before writing so we never clobber the wrong cell.\nconst replacements = [\n  [\"75\u00f78=9, 3\", \"40\u00f75=8, 0\"],\n  [\"37\u00f78=4, 5\", \"62\u00f74=15, 2\"],\n  [\"71\u00f74=17, 3\", \"55\u00f74=13, 3\"],\n  [\"60\u00f72=30, 0\", \"18\u00f74=4, 2\"],\n  [\"85\u00f72=42, 1\", \"15\u00f72=7, 1\"],\n  [\"29\u00f73=9, 2\", \"72\u00f75=14, 2\"],\n  [\"53\u00f76=8, 5\", \"36\u00f72=18, 0\"],\n  [\"58\u00f77=8, 2\", \"11\u00f73=3, 2\"],\n  [\"69\u00f76=11, 3\", \"48\u00f73=16, 0\"],\n  [\"38\u00f75=7, 3\", \"29\u00f78=3, 5\"],\n  [\"73\u00f75=14, 3\", \"85\u00f75=17, 0\"],\n  [\"68\u00f72=34, 0\", \"98\u00f76=16, 2\"],\n  [\"33\u00f77=4, 5\", \"91\u00f73=30, 1\"],\n  [\"55\u00f75=11, 0\", \"29\u00f72=14, 1\"],\n  [\"99\u00f77=14, 1\", \"35\u00f72=17, 1\"],\n  [\"95\u00f74=23, 3\", \"30\u00f78=3, 6\"],\n  [\"92\u00f72=46, 0\", \"58\u00f75=11, 3\"],\n  [\"51\u00f77=7, 2\", \"84\u00f75=16, 4\"],\n  [\"50\u00f77=7, 1\", \"19\u00f77=2, 5\"],\n  [\"59\u00f72=29, 1\", \"29\u00f77=4, 1\"],\n  [\"48\u00f74=12, 0\", \"73\u00f74=18, 1\"],\n  [\"73\u00f76=12, 1\", \"52\u00f73=17, 1\"],\n  [\"95\u00f77=13, 4\", \"81\u00f76=13, 3\"],\n  [\"64\u00f78=8, 0\", \"12\u00f72=6, 0\"],\n  [\"33\u00f79=3, 6\", \"59\u00f72=29, 1\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Load every cell's current value up front.\nconst cellsByRow = [];\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n  cellsByRow.push(row.cells);\n}\nawait context.sync();\n\nfor (const cells of cellsByRow) {\n  for (const cell of cells.items) {\n    cell.load(\"value\");\n  }\n}\nawait context.sync();\n\n// Walk cells in row-major order, consuming the next replacement whenever\n// a cell's current text matches the expected \"old\" value.\nlet idx = 0;\nfor (const cells of cellsByRow) {\n  for (const cell of cells.items) {\n    if (idx >= replacements.length) break;\n    const [oldText, newText] = replacements[idx];\n    if (cell.value === oldText) {\n      cell.value = newText;\n      idx++;\n    }\n  }\n}\nawait context.sync();\n\nif (idx !== replacements.length) {\n  throw new Error(\n    `Only applied ${idx} of ${replacements.length} replacements; ` +\n      `document text did not match the expected pattern.`\n  );\n}\n", "ps1": "# Replace the 25 two-digit-division answers in the single table, in\n# document (row-major) order. Only every 4th row (1, 5, 9, 13, 17) holds\n# text; the rows between are blank spacer rows. Each entry is\n# (oldText, newText); the old text is asserted before writing so we never\n# clobber the wrong cell.\n$replacements = @(\n  @(\"75\u00f78=9, 3\", \"40\u00f75=8, 0\"),\n  @(\"37\u00f78=4, 5\", \"62\u00f74=15, 2\"),\n  @(\"71\u00f74=17, 3\", \"55\u00f74=13, 3\"),\n  @(\"60\u00f72=30, 0\", \"18\u00f74=4, 2\"),\n  @(\"85\u00f72=42, 1\", \"15\u00f72=7, 1\"),\n  @(\"29\u00f73=9, 2\", \"72\u00f75=14, 2\"),\n  @(\"53\u00f76=8, 5\", \"36\u00f72=18, 0\"),\n  @(\"58\u00f77=8, 2\", \"11\u00f73=3, 2\"),\n  @(\"69\u00f76=11, 3\", \"48\u00f73=16, 0\"),\n  @(\"38\u00f75=7, 3\", \"29\u00f78=3, 5\"),\n  @(\"73\u00f75=14, 3\", \"85\u00f75=17, 0\"),\n  @(\"68\u00f72=34, 0\", \"98\u00f76=16, 2\"),\n  @(\"33\u00f77=4, 5\", \"91\u00f73=30, 1\"),\n  @(\"55\u00f75=11, 0\", \"29\u00f72=14, 1\"),\n  @(\"99\u00f77=14, 1\", \"35\u00f72=17, 1\"),\n  @(\"95\u00f74=23, 3\", \"30\u00f78=3, 6\"),\n  @(\"92\u00f72=46, 0\", \"58\u00f75=11, 3\"),\n  @(\"51\u00f77=7, 2\", \"84\u00f75=16, 4\"),\n  @(\"50\u00f77=7, 1\", \"19\u00f77=2, 5\"),\n  @(\"59\u00f72=29, 1\", \"29\u00f77=4, 1\"),\n  @(\"48\u00f74=12, 0\", \"73\u00f74=18, 1\"),\n  @(\"73\u00f76=12, 1\", \"52\u00f73=17, 1\"),\n  @(\"95\u00f77=13, 4\", \"81\u00f76=13, 3\"),\n  @(\"64\u00f78=8, 0\", \"12\u00f72=6, 0\"),\n  @(\"33\u00f79=3, 6\", \"59\u00f72=29, 1\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$idx = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n  for ($c = 1; $c -le $t.Columns.Count; $c++) {\n    if ($idx -ge $replacements.Count) { break }\n    $cell = $t.Cell($r, $c)\n    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    $pair = $replacements[$idx]\n    if ($current -eq $pair[0]) {\n      $cell.Range.Text = $pair[1]\n      $idx++\n    }\n  }\n}\n\nif ($idx -ne $replacements.Count) {\n  throw \"Only applied $idx of $($replacements.Count) replacements; document text did not match the expected pattern.\"\n}\n"}
